$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 5 data (episode 9/10 results) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Shelly"
$ws.Range("C5").Value = "Jana"
$ws.Range("D5").Value = "Pia"
$ws.Range("E5").Value = "Edda"
$ws.Range("F5").Value = "Julia"
$ws.Range("G5").Value = "Tais"
$ws.Range("H5").Value = "Sina"
$ws.Range("I5").Value = "Lina"
$ws.Range("J5").Value = "Afra"
$ws.Range("K5").Value = "Lisa-Marie"
$ws.Range("L5").Value = 3

# --- Highlight cells: green (correct) fill ---
$greenCells = @("B2", "G3", "B4", "G4", "B5", "G5")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5287936
}

# --- Highlight cells: red (wrong) fill ---
$redCells = @("E2", "G2", "K2", "B3", "I3", "H4")
foreach ($addr in $redCells) {
    $ws.Range($addr).Interior.Color = 255
}

# --- Update selection / active cell ---
$ws.Range("H4").Select()
